$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 640, shifting existing rows 640:749 down to 641:750
$ws.Rows.Item(640).Insert()

# Populate the newly inserted row 640 with its data
$ws.Range("A640").Value = 10
$ws.Range("B640").Value = "Vega Modelo de Temuco"
$ws.Range("C640").Value = "La Araucanía"
$ws.Range("D640").Value = 44951
$ws.Range("E640").Value = 9
$ws.Range("F640").Value = 100112045
$ws.Range("G640").Value = "Zapallo"
$ws.Range("H640").Value = "Camote"
$ws.Range("I640").Value = "1a (cosecha)"
$ws.Range("J640").Value = 450
$ws.Range("K640").Value = 700
$ws.Range("L640").Value = 700
$ws.Range("M640").Value = 700
$ws.Range("N640").Value = "$/kilo (volumen en unidades)"
$ws.Range("O640").Value = "Región del Maule"
$ws.Range("P640").Value = 700
$ws.Range("Q640").Value = 1
$ws.Range("R640").Value = "Hortaliza"
